$d = $word.ActiveDocument

# --- Change 1 -------------------------------------------------------
# Merge the two runs "THU Sep 11" + " 11:27:50 IST 2018" into a single
# run "THU Sep 11 11:27:50 IST 2018" (a Find/Replace on the exact
# concatenated text merges the split runs into one while it rewrites
# the text).
$null = $d.Content.Find.Execute("THU Sep 11 11:27:50 IST 2018", $true, $false, $false, $false, $false, $true, 1, $false, "THU Sep 11 11:27:50 IST 2018", 2)

# --- Change 2 -------------------------------------------------------
# Append a brand-new "chick in" purchase-details block (dated
# TUE Sep 18) right after the last filled-in record (the final
# "Amount balance ... - 530.0" paragraph), keeping all the blank
# trailing paragraphs that already existed at the end of the story.
$count = $d.Paragraphs.Count
$targetIndex = -1
for ($i = $count; $i -ge 1; $i--) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text.Trim().Length -gt 0) {
        $targetIndex = $i
        break
    }
}
$lastPara = $d.Paragraphs.Item($targetIndex)
$insertPoint = $d.Range($lastPara.Range.End, $lastPara.Range.End)

$fontRpr = '<w:rFonts w:ascii="Courier New" w:hAnsi="Courier New" w:cs="Courier New"/>'
$fontRprBold = '<w:rFonts w:ascii="Courier New" w:hAnsi="Courier New" w:cs="Courier New"/><w:b/>'
$fontRprRed = '<w:rFonts w:ascii="Courier New" w:hAnsi="Courier New" w:cs="Courier New"/><w:color w:val="FF0000"/>'

$newBlockXml = @"
<w:p><w:pPr><w:pStyle w:val="PlainText"/><w:rPr>$fontRprBold</w:rPr></w:pPr></w:p>
<w:p><w:pPr><w:pStyle w:val="PlainText"/><w:rPr>$fontRpr</w:rPr></w:pPr>
  <w:r><w:rPr>$fontRpr</w:rPr><w:t>TUE Sep 18</w:t></w:r>
  <w:r><w:rPr>$fontRpr</w:rPr><w:t xml:space="preserve"> 11:30:25 IST 2018</w:t></w:r>
</w:p>
<w:p><w:pPr><w:pStyle w:val="PlainText"/><w:rPr>$fontRpr</w:rPr></w:pPr>
  <w:r><w:rPr>$fontRpr</w:rPr><w:t>Person Name</w:t></w:r>
  <w:r><w:rPr>$fontRpr</w:rPr><w:tab/></w:r>
  <w:r><w:rPr>$fontRpr</w:rPr><w:tab/></w:r>
  <w:r><w:rPr>$fontRpr</w:rPr><w:tab/></w:r>
  <w:r><w:rPr>$fontRpr</w:rPr><w:tab/><w:t>- NAGAMMA</w:t></w:r>
</w:p>
<w:p><w:pPr><w:pStyle w:val="PlainText"/><w:rPr>$fontRpr</w:rPr></w:pPr>
  <w:r><w:rPr>$fontRpr</w:rPr><w:t>Bill number</w:t></w:r>
  <w:r><w:rPr>$fontRpr</w:rPr><w:tab/></w:r>
  <w:r><w:rPr>$fontRpr</w:rPr><w:tab/></w:r>
  <w:r><w:rPr>$fontRpr</w:rPr><w:tab/></w:r>
  <w:r><w:rPr>$fontRpr</w:rPr><w:tab/><w:t>- 7890</w:t></w:r>
</w:p>
<w:p><w:pPr><w:pStyle w:val="PlainText"/><w:rPr>$fontRpr</w:rPr></w:pPr>
  <w:r><w:rPr>$fontRpr</w:rPr><w:t>---------------------------------------------------------------</w:t></w:r>
</w:p>
<w:p><w:pPr><w:pStyle w:val="PlainText"/><w:rPr>$fontRpr</w:rPr></w:pPr>
  <w:r><w:rPr>$fontRpr</w:rPr><w:t>Item Name</w:t></w:r>
  <w:r><w:rPr>$fontRpr</w:rPr><w:tab/></w:r>
  <w:r><w:rPr>$fontRpr</w:rPr><w:tab/></w:r>
  <w:r><w:rPr>$fontRpr</w:rPr><w:tab/></w:r>
  <w:r><w:rPr>$fontRpr</w:rPr><w:tab/><w:t>- BEET</w:t></w:r>
</w:p>
<w:p><w:pPr><w:pStyle w:val="PlainText"/><w:rPr>$fontRpr</w:rPr></w:pPr>
  <w:r><w:rPr>$fontRpr</w:rPr><w:t>Number of Pockets</w:t></w:r>
  <w:r><w:rPr>$fontRpr</w:rPr><w:tab/></w:r>
  <w:r><w:rPr>$fontRpr</w:rPr><w:tab/></w:r>
  <w:r><w:rPr>$fontRpr</w:rPr><w:tab/><w:t>- 1</w:t></w:r>
</w:p>
<w:p><w:pPr><w:pStyle w:val="PlainText"/><w:rPr>$fontRpr</w:rPr></w:pPr>
  <w:r><w:rPr>$fontRpr</w:rPr><w:t>Number of KGs</w:t></w:r>
  <w:r><w:rPr>$fontRpr</w:rPr><w:tab/></w:r>
  <w:r><w:rPr>$fontRpr</w:rPr><w:tab/></w:r>
  <w:r><w:rPr>$fontRpr</w:rPr><w:tab/><w:t>- 55</w:t></w:r>
</w:p>
<w:p><w:pPr><w:pStyle w:val="PlainText"/><w:rPr>$fontRpr</w:rPr></w:pPr>
  <w:r><w:rPr>$fontRpr</w:rPr><w:t>Rate</w:t></w:r>
  <w:r><w:rPr>$fontRpr</w:rPr><w:tab/></w:r>
  <w:r><w:rPr>$fontRpr</w:rPr><w:tab/></w:r>
  <w:r><w:rPr>$fontRpr</w:rPr><w:tab/></w:r>
  <w:r><w:rPr>$fontRpr</w:rPr><w:tab/></w:r>
  <w:r><w:rPr>$fontRpr</w:rPr><w:tab/><w:t>- 10</w:t></w:r>
</w:p>
<w:p><w:pPr><w:pStyle w:val="PlainText"/><w:rPr>$fontRpr</w:rPr></w:pPr>
  <w:r><w:rPr>$fontRpr</w:rPr><w:t>Total Price</w:t></w:r>
  <w:r><w:rPr>$fontRpr</w:rPr><w:tab/></w:r>
  <w:r><w:rPr>$fontRpr</w:rPr><w:tab/></w:r>
  <w:r><w:rPr>$fontRpr</w:rPr><w:tab/></w:r>
  <w:r><w:rPr>$fontRpr</w:rPr><w:tab/><w:t>- 550.0</w:t></w:r>
</w:p>
<w:p><w:pPr><w:pStyle w:val="PlainText"/><w:rPr>$fontRprRed</w:rPr></w:pPr>
  <w:r><w:rPr>$fontRprRed</w:rPr><w:t>Amount Received</w:t></w:r>
  <w:r><w:rPr>$fontRprRed</w:rPr><w:tab/></w:r>
  <w:r><w:rPr>$fontRprRed</w:rPr><w:tab/></w:r>
  <w:r><w:rPr>$fontRprRed</w:rPr><w:tab/><w:t>- 530</w:t></w:r>
</w:p>
<w:p><w:pPr><w:pStyle w:val="PlainText"/><w:rPr>$fontRprBold</w:rPr></w:pPr>
  <w:r><w:rPr>$fontRprBold</w:rPr><w:t>Amount balance</w:t></w:r>
  <w:r><w:rPr>$fontRprBold</w:rPr><w:tab/></w:r>
  <w:r><w:rPr>$fontRprBold</w:rPr><w:tab/></w:r>
  <w:r><w:rPr>$fontRprBold</w:rPr><w:tab/><w:t>- 550.0</w:t></w:r>
</w:p>
<w:p><w:pPr><w:pStyle w:val="PlainText"/><w:rPr>$fontRprBold</w:rPr></w:pPr>
  <w:r><w:rPr>$fontRprBold</w:rPr><w:t>Amount Received mode</w:t></w:r>
  <w:r><w:rPr>$fontRprBold</w:rPr><w:tab/></w:r>
  <w:r><w:rPr>$fontRprBold</w:rPr><w:tab/><w:t>- CASH</w:t></w:r>
</w:p>
<w:p><w:pPr><w:pStyle w:val="PlainText"/><w:rPr>$fontRpr</w:rPr></w:pPr></w:p>
<w:p><w:pPr><w:pStyle w:val="PlainText"/><w:rPr>$fontRprBold</w:rPr></w:pPr></w:p>
"@

$null = $insertPoint.InsertXML($newBlockXml)
